# Bulk OT Booking Import.xlsx - remove the ApplicationTypeName column and
# add a couple of formatted "helper" cells used to pin down a date format
# and a time format (mirrors the author's edit to the template).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "ApplicationTypeName" header/column (col F) - it's no longer
# part of the import template. Do this before touching anything further
# right, so later writes land on their final (post-delete) addresses.
$ws.Range("F1").EntireColumn.Delete()

# Re-assert the (now best-fit) widths for the remaining data columns.
$ws.Columns("A:A").ColumnWidth = 5.833333333333333
$ws.Columns("B:B").ColumnWidth = 6.666666666666667
$ws.Columns("C:C").ColumnWidth = 8.666666666666666
$ws.Columns("D:D").ColumnWidth = 7.666666666666667
$ws.Columns("E:E").ColumnWidth = 6.666666666666667

# Stamp a date-formatted cell and a time-formatted cell on row 2 so the
# template carries the expected number formats for OTDate / StartTime.
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("J2").NumberFormat = "mm:ss.0"

# Match the author's final selection in the sheet.
[void]$ws.Range("F3").Select()
